$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D5").Value = "2016-02-17 03:27:22"
$wsZh.Range("G5").Value = "2016-02-17 03:28:05"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D5").Value = "2016-02-17 03:27:33"
$wsDe.Range("G5").Value = "2016-02-17 03:28:22"
